# Update robot-wiring workbook: replace the old left/right/back sonar trig+echo
# rows with a new sonarOffset-style model that lists front/left/right/back
# trig/echo pins as their own rows at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the old sonar pin-number + description cells from rows 10-15
# (pin index in column A stays as-is, it's just the wiring description that
# is being removed/relocated).
$ws.Range("B10:C10").Value = $null
$ws.Range("B11:C11").Value = $null
$ws.Range("B12:C12").Value = $null
$ws.Range("B13:C13").Value = $null
$ws.Range("B14:C14").Value = $null
$ws.Range("B15:C15").Value = $null

# Append the new sonarOffset rows (22-30) describing front/left/right/back
# trig and echo wiring.
$ws.Range("A24").Value = 22
$ws.Range("C24").Value = "left trig"

$ws.Range("A25").Value = 23
$ws.Range("C25").Value = "left echo "

$ws.Range("A26").Value = 24
$ws.Range("C26").Value = "right trig"

$ws.Range("A27").Value = 25
$ws.Range("C27").Value = "right echo "

$ws.Range("A28").Value = 26
$ws.Range("C28").Value = "back trig "

$ws.Range("A29").Value = 27
$ws.Range("C29").Value = "back echo "

$ws.Range("A30").Value = 28
$ws.Range("C30").Value = "front trig"

$ws.Range("A31").Value = 29
$ws.Range("C31").Value = "front echo"

$ws.Range("A32").Value = 30

# Update the view: clear the previous scroll position / selection
# (was topLeftCell A7 / selection D16) to a plain selection on D8.
$ws.Range("D8").Select()
